$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the bold/centered/bordered header format (currently on B1:F1) before
# any values move around, then stamp it across the future header row A1:E1.
$ws.Range("B1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

# Shift the header row one column to the left (old B1:F1 -> new A1:E1).
# While doing so, also fix the MODEL_CONDITION -> MODELCONDITION typo.
$ws.Range("A1").Value = "EL_Astral50"
$ws.Range("B1").Value = "FNRATE_PHYLONET"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODELCONDITION"
$ws.Range("E1").Value = "GENE"

# Shift row 2 one column to the left (old A2:F2 -> new A2:E2), dropping the
# old A2 taxon-count value and keeping the rest.
$ws.Range("A2").Value = 170
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "simulated_50genes_weakILS"
$ws.Range("E2").Value = 1

# Shift row 3 one column to the left the same way.
$ws.Range("A3").Value = 170
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "simulated_50genes_weakILS"
$ws.Range("E3").Value = 15

# The old column F is no longer used.
$ws.Range("F1:F3").Clear()

# The data rows (2 and 3) should carry no special formatting; the old A2/A3
# cells had the bordered header style which must not carry over.
$ws.Range("A2:E3").ClearFormats()
